$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.644.07"
$ws.Range("E2").Value = "  -1.56%  "
$ws.Range("D3").Value = "1.878.74"
$ws.Range("E3").Value = "  -1.35%  "
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "330.16"
$ws.Range("E5").Value = "  +1.02%  "
$ws.Range("E6").Value = "  -0.09%  "
$ws.Range("D7").Value = "0.4717"
$ws.Range("E7").Value = "  +1.83%  "
$ws.Range("D8").Value = "0.3980"
$ws.Range("E8").Value = "  +0.17%  "
$ws.Range("D9").Value = "48.44"
$ws.Range("E9").Value = "  -7.11%  "
$ws.Range("D10").Value = "0.08052"
$ws.Range("E10").Value = "  -3.22%  "
$ws.Range("D11").Value = "1.025"
$ws.Range("E11").Value = "  -1.64%  "
$ws.Range("D12").Value = "21.83"
$ws.Range("E12").Value = "  -0.13%  "
$ws.Range("D13").Value = "1.885.35"
$ws.Range("E13").Value = "  -0.92%  "
$ws.Range("D14").Value = "5.960"
$ws.Range("E14").Value = "  -0.93%  "
$ws.Range("D15").Value = "7.195"
$ws.Range("E15").Value = "  -2.44%  "
$ws.Range("E16").Value = "  -0.02%  "
$ws.Range("D17").Value = "86.93"
$ws.Range("E17").Value = "  -2.66%  "
$ws.Range("E18").Value = "  -2.32%  "
$ws.Range("D19").Value = "0.06572"
$ws.Range("E19").Value = "  -0.15%  "
$ws.Range("D20").Value = "17.22"
$ws.Range("E20").Value = "  -3.37%  "
$ws.Range("E21").Value = "  -0.07%  "
$ws.Range("D22").Value = "27.674.29"
$ws.Range("E22").Value = "  -1.43%  "
$ws.Range("D23").Value = "5.508"
$ws.Range("E23").Value = "  -3.44%  "
$ws.Range("D24").Value = "10.99"
$ws.Range("E24").Value = "  -1.39%  "
$ws.Range("E25").Value = "  -0.22%  "
$ws.Range("D26").Value = "2.104.74"
$ws.Range("E26").Value = "  -1.23%  "
$ws.Range("D27").Value = "154.46"
$ws.Range("E27").Value = "  +0.30%  "
$ws.Range("D28").Value = "20.27"
$ws.Range("E28").Value = "  +1.32%  "
$ws.Range("D29").Value = "2.099"
$ws.Range("E29").Value = "  -1.19%  "
$ws.Range("D30").Value = "5.604"
$ws.Range("E30").Value = "  -2.09%  "
$ws.Range("D31").Value = "122.46"
$ws.Range("E31").Value = "  -1.57%  "
$ws.Range("D32").Value = "0.09507"
$ws.Range("E32").Value = "  -1.06%  "
$ws.Range("D33").Value = "0.9545"
$ws.Range("E33").Value = "  -1.39%  "
$ws.Range("D34").Value = "1.474"
$ws.Range("E34").Value = "  +0.36%  "
$ws.Range("D35").Value = "3.615"
$ws.Range("E35").Value = "  -0.23%  "
$ws.Range("D36").Value = "5.306"
$ws.Range("E36").Value = "  -3.67%  "
$ws.Range("D37").Value = "0.06109"
$ws.Range("E37").Value = "  -0.52%  "
$ws.Range("D38").Value = "0.02253"
$ws.Range("E38").Value = "  -1.47%  "
$ws.Range("E39").Value = "  -3.20%  "
$ws.Range("D40").Value = "8.233"
$ws.Range("E40").Value = "  -5.53%  "
$ws.Range("D41").Value = "0.5995"
$ws.Range("E41").Value = "  -2.13%  "
$ws.Range("E42").Value = "  -0.06%  "
$ws.Range("D43").Value = "0.1898"
$ws.Range("E43").Value = "  -0.16%  "
$ws.Range("D44").Value = "10.34"
$ws.Range("E44").Value = "  -4.77%  "
$ws.Range("D45").Value = "0.5699"
$ws.Range("E45").Value = "  -2.60%  "
$ws.Range("E46").Value = "  -5.13%  "
$ws.Range("E47").Value = "  -4.46%  "
$ws.Range("D48").Value = "3.411"
$ws.Range("E48").Value = "  -0.60%  "
$ws.Range("D49").Value = "1.941"
$ws.Range("E49").Value = "  -3.31%  "
$ws.Range("D50").Value = "0.06823"
$ws.Range("E50").Value = "  -1.10%  "
$ws.Range("D51").Value = "110.11"
$ws.Range("E51").Value = "  -0.60%  "
